# Update the "formulas" sheet: re-pad the existing vpcNNN bundle labels to
# 4-digit numbers (vpc0001 ... vpc0060) and extend the bundle list down to
# vpc0108, adding rows 63-110 with the same alternating row style used by
# the existing table (even row -> style of row 2/4/6..., odd row -> style
# of row 3/5/7...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formulas")

# Style references already present on the sheet for the alternating
# "ws_even"/"ws_odd" look used by the bundle column.
$evenStyle = $ws.Cells.Item(4, 1).Style   # style used on even rows (s="1")
$oddStyle  = $ws.Cells.Item(3, 1).Style   # style used on odd rows  (s="2")

for ($row = 3; $row -le 110; $row++) {
    $num = $row - 2
    $label = "vpc" + $num.ToString().PadLeft(4, '0')
    $ws.Cells.Item($row, 3).Value = $label

    if ($row -gt 62) {
        if (($row % 2) -eq 0) {
            $style = $evenStyle
        } else {
            $style = $oddStyle
        }
        for ($col = 1; $col -le 4; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $cell.Style = $style
            # Re-assert alignment explicitly so the engine reuses the exact
            # existing style index instead of minting a visually-identical
            # but distinct one (its Style copy alone drops the alignment).
            $cell.HorizontalAlignment = -4108
            $cell.VerticalAlignment = -4108
        }
    }
}

Write-Host "done"
